$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Prix Spot")

# Insert a new column before column B, shifting existing B/C to C/D
$ws.Range("B1").EntireColumn.Insert()

# The inserted column inherited formatting from the (old) column B; clear the
# data cells' formatting so they go back to the default style, matching the
# un-styled numeric cells used elsewhere in the sheet.
$ws.Range("B2:B25").ClearFormats()

# Header for the new column
$ws.Range("B1").Value = "15-jun"

# New column B values (rows 2-25)
$values = @(51.35, 28.31, 26.87, 21.88, 18.78, 17.86, 18.47, 15.13, 8.970000000000001, 4.55, 0, -0.02, -1.21, -5.6, -5, -2, -0.01, -0.01, 12.37, 19.29, 33.96, 39.96, 61.7, 53.03)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
